# Update the log write mode: refresh simulated run metrics
# (run_time, num_deaths, max_er, and the per-iteration convergence
# series) for gr100_08_simulated/details with the latest re-run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.7444260120391846
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 5818.209725119707
$ws.Range("I2").Value = 0.1612467734209603
$ws.Range("J2").Value = 0.1432357924831252
$ws.Range("K2").Value = 0.1409509607301358
$ws.Range("L2").Value = 0.1373888944078355
$ws.Range("M2").Value = 0.1358179742418151
$ws.Range("N2").Value = 0.1294153942518461
$ws.Range("O2").Value = 0.1294153942518461
$ws.Range("P2").Value = 0.1294153942518461
$ws.Range("Q2").Value = 0.1294153942518461
$ws.Range("R2").Value = 0.1294153942518461
$ws.Range("S2").Value = 0.1294153942518461
$ws.Range("T2").Value = 0.1294153942518461
$ws.Range("U2").Value = 0.1294153942518461
$ws.Range("V2").Value = 0.1294153942518461
$ws.Range("W2").Value = 0.1294153942518461
$ws.Range("X2").Value = 0.1294153942518461
$ws.Range("Y2").Value = 0.1294153942518461

# Row 3
$ws.Range("C3").Value = 1.078359365463257
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 5549.580757867569
$ws.Range("I3").Value = 0.1649773832155262
$ws.Range("J3").Value = 0.1509042501999619
$ws.Range("K3").Value = 0.1453683314825268
$ws.Range("L3").Value = 0.1319989883509791
$ws.Range("M3").Value = 0.1319989883509791
$ws.Range("N3").Value = 0.1318148275497462
$ws.Range("O3").Value = 0.123428847305307
$ws.Range("P3").Value = 0.1222815401950872
$ws.Range("Q3").Value = 0.1198868791837298
$ws.Range("R3").Value = 0.1198868791837298
$ws.Range("S3").Value = 0.1194798072655575
$ws.Range("T3").Value = 0.1185372427381384
$ws.Range("U3").Value = 0.1176361161952569
$ws.Range("V3").Value = 0.1172032778254521
$ws.Range("W3").Value = 0.116178962141668
$ws.Range("X3").Value = 0.116178962141668
$ws.Range("Y3").Value = 0.116178962141668

# Row 4
$ws.Range("C4").Value = 0.9218385219573975
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 5522.031051209871
$ws.Range("I4").Value = 0.165726156694249
$ws.Range("J4").Value = 0.1556224032727453
$ws.Range("K4").Value = 0.1407074144722224
$ws.Range("L4").Value = 0.1332354792257785
$ws.Range("M4").Value = 0.1304841187981101
$ws.Range("N4").Value = 0.1258712414116389
$ws.Range("O4").Value = 0.1258712414116389
$ws.Range("P4").Value = 0.1240950807058515
$ws.Range("Q4").Value = 0.1240156270650122
$ws.Range("R4").Value = 0.123898973301136
$ws.Range("S4").Value = 0.123898973301136
$ws.Range("T4").Value = 0.1236891168188872
$ws.Range("U4").Value = 0.1236891168188872
$ws.Range("V4").Value = 0.1236891168188872
$ws.Range("W4").Value = 0.1236891168188872
$ws.Range("X4").Value = 0.1236419308228045
$ws.Range("Y4").Value = 0.1236419308228045

# Row 5
$ws.Range("C5").Value = 0.8437650203704834
$ws.Range("E5").Value = 5528.375565826568
$ws.Range("K5").Value = 0.1355208441974171
$ws.Range("L5").Value = 0.1328800243010802
$ws.Range("M5").Value = 0.1289915434127319
$ws.Range("N5").Value = 0.1220340798249276
$ws.Range("O5").Value = 0.1220340798249276
$ws.Range("P5").Value = 0.1197445999028133
$ws.Range("Q5").Value = 0.1188657219056261
$ws.Range("R5").Value = 0.1188657219056261
$ws.Range("S5").Value = 0.1184552727020081
$ws.Range("T5").Value = 0.1173854298325454
$ws.Range("U5").Value = 0.11674686560849
$ws.Range("V5").Value = 0.1167089826358899
$ws.Range("W5").Value = 0.115765605571668
$ws.Range("X5").Value = 0.115765605571668
$ws.Range("Y5").Value = 0.115765605571668

# Row 6
$ws.Range("C6").Value = 0.8593864440917969
$ws.Range("E6").Value = 5777.664946093589
$ws.Range("I6").Value = 0.1556455539133602
$ws.Range("J6").Value = 0.1513237459422206
$ws.Range("K6").Value = 0.1436460872980684
$ws.Range("L6").Value = 0.1351273369360836
$ws.Range("M6").Value = 0.1351273369360836
$ws.Range("N6").Value = 0.1338060286462254
$ws.Range("O6").Value = 0.1247641134643513
$ws.Range("P6").Value = 0.1245414845932979
$ws.Range("Q6").Value = 0.123262470692648
$ws.Range("R6").Value = 0.1229008870636346
$ws.Range("S6").Value = 0.1229008870636346
$ws.Range("T6").Value = 0.1219306882078127
$ws.Range("U6").Value = 0.1212731419103084
$ws.Range("V6").Value = 0.1210699655274099
$ws.Range("W6").Value = 0.1209526527960076
$ws.Range("X6").Value = 0.1208039751611709
$ws.Range("Y6").Value = 0.1206250476821362

# Row 7
$ws.Range("C7").Value = 0.843724250793457
$ws.Range("E7").Value = 5635.290170927763
$ws.Range("I7").Value = 0.16053514371711
$ws.Range("J7").Value = 0.1418392959622381
$ws.Range("K7").Value = 0.137046920160777
$ws.Range("L7").Value = 0.1351636666743501
$ws.Range("M7").Value = 0.1324155919062014
$ws.Range("N7").Value = 0.1313625134089997
$ws.Range("O7").Value = 0.1278701804191927
$ws.Range("P7").Value = 0.1278701804191927
$ws.Range("Q7").Value = 0.1268864056362424
$ws.Range("R7").Value = 0.1265627450762096
$ws.Range("S7").Value = 0.1259587040413774
$ws.Range("T7").Value = 0.1259587040413774
$ws.Range("U7").Value = 0.1259587040413774
$ws.Range("V7").Value = 0.1259587040413774
$ws.Range("W7").Value = 0.1258497109342643
$ws.Range("X7").Value = 0.1258497109342643
$ws.Range("Y7").Value = 0.1258497109342643

# Row 8
$ws.Range("C8").Value = 0.8594014644622803
$ws.Range("E8").Value = 5561.121790263073
$ws.Range("I8").Value = 0.165726156694249
$ws.Range("J8").Value = 0.1426082550982998
$ws.Range("K8").Value = 0.1415162970344212
$ws.Range("L8").Value = 0.1397365322063736
$ws.Range("M8").Value = 0.1326066613395584
$ws.Range("N8").Value = 0.124597596163596
$ws.Range("O8").Value = 0.124597596163596
$ws.Range("P8").Value = 0.1236637241515554
$ws.Range("Q8").Value = 0.1193583876122125
$ws.Range("R8").Value = 0.1193583876122125
$ws.Range("S8").Value = 0.1192216178901757
$ws.Range("T8").Value = 0.117335653469729
$ws.Range("U8").Value = 0.117335653469729
$ws.Range("V8").Value = 0.116749629832031
$ws.Range("W8").Value = 0.116749629832031
$ws.Range("X8").Value = 0.116749629832031
$ws.Range("Y8").Value = 0.1164039335333932

# Row 9
$ws.Range("C9").Value = 0.7655999660491943
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 5566.681291594658
$ws.Range("J9").Value = 0.1492580251219544
$ws.Range("K9").Value = 0.1371595902975441
$ws.Range("L9").Value = 0.1300566062269913
$ws.Range("M9").Value = 0.1252352845453791
$ws.Range("N9").Value = 0.124485728366017
$ws.Range("O9").Value = 0.124485728366017
$ws.Range("P9").Value = 0.1234453567825319
$ws.Range("Q9").Value = 0.121982845543144
$ws.Range("R9").Value = 0.1216208724505044
$ws.Range("S9").Value = 0.1202929609455025
$ws.Range("T9").Value = 0.119218317748401
$ws.Range("U9").Value = 0.1188849911772946
$ws.Range("V9").Value = 0.1187661615536152
$ws.Range("W9").Value = 0.1170113783455089
$ws.Range("X9").Value = 0.1168233446628406
$ws.Range("Y9").Value = 0.1165123058790381

# Row 10
$ws.Range("C10").Value = 0.7969136238098145
$ws.Range("E10").Value = 5523.824040974669
$ws.Range("I10").Value = 0.165726156694249
$ws.Range("J10").Value = 0.1538881665131734
$ws.Range("K10").Value = 0.1472399347538996
$ws.Range("L10").Value = 0.1387525116533764
$ws.Range("M10").Value = 0.136538133377005
$ws.Range("N10").Value = 0.1332828548476813
$ws.Range("O10").Value = 0.1284266502415519
$ws.Range("P10").Value = 0.1252981294095448
$ws.Range("Q10").Value = 0.1250570663618751
$ws.Range("R10").Value = 0.1250570663618751
$ws.Range("S10").Value = 0.124873683282672
$ws.Range("T10").Value = 0.124873683282672
$ws.Range("U10").Value = 0.1239702349983815
$ws.Range("V10").Value = 0.1239702349983815
$ws.Range("W10").Value = 0.1239702349983815
$ws.Range("X10").Value = 0.1239702349983815
$ws.Range("Y10").Value = 0.1236768818903444

# Row 11
$ws.Range("C11").Value = 0.7499873638153076
$ws.Range("E11").Value = 5384.813232915303
$ws.Range("J11").Value = 0.1556224032727453
$ws.Range("K11").Value = 0.1419771142015106
$ws.Range("L11").Value = 0.1363032853960092
$ws.Range("M11").Value = 0.1289843925104508
$ws.Range("N11").Value = 0.1226246923756752
$ws.Range("O11").Value = 0.1224755686838483
$ws.Range("P11").Value = 0.1208451679664848
$ws.Range("Q11").Value = 0.1187675514144697
$ws.Range("R11").Value = 0.1174163406949813
$ws.Range("S11").Value = 0.1161612618409982
$ws.Range("T11").Value = 0.1155934138523988
$ws.Range("U11").Value = 0.1150157652565395
$ws.Range("V11").Value = 0.1129774923165696
$ws.Range("W11").Value = 0.1129774923165696
$ws.Range("X11").Value = 0.1129671195500059
